$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph entirely (it was right
#    after the "Play Braccio di Ferro Slot Game for Free" heading).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Meta description")) {
        $deleted = $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Braccio di Ferro Slot Game for
#    Free" right before the "Prompt: ..." paragraph near the end of the
#    document.
# ---------------------------------------------------------------------
$promptPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Prompt:")) {
        $promptPara = $p
        break
    }
}

if ($promptPara -ne $null) {
    $prevPara = $promptPara.Previous()
    # Insertion point just before the paragraph mark that ends $prevPara,
    # i.e. right at the boundary between $prevPara and $promptPara.
    $insertAt = $prevPara.Range.End - 1
    $insertRange = $d.Range($insertAt, $insertAt)

    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Braccio di Ferro Slot Game for Free</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $inserted = $insertRange.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 3) Replace the old "Prompt: ..." text with the new meta-description
#    style sentence (keeping the run's italic formatting intact).
# ---------------------------------------------------------------------
$old = 'Prompt: Create a feature image for the Braccio di Ferro slot game that captures its adventurous and playful spirit. The image should be in cartoon style and include a happy Maya warrior with glasses. The Maya warrior should be depicted engaging with the game, possibly spinning the reels or collecting flying fish in a bottle. The background should be a small harbor with boats and seagulls in the sky. Use bright colors and playful elements to reflect the fun and excitement of the game. Make sure to include the name of the game, "Braccio di Ferro", in an eye-catching font.'
$new = 'Explore the gameplay features, bonuses, graphics, wagering options, and RTP rate of Braccio di Ferro slot game and play it for free.'
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
